# "Added last minute updates"
# The doc's first paragraph (the hidden **ID__...__ID** bookmark marker)
# gets re-pointed at the AF_PGI_5316_506 topic id, picks up the same
# "no-line" paragraph border spacing already used by the two FAR/AFFARS
# body paragraphs below it, and its left indent grows from 120 -> 225
# twips (6pt -> 11.25pt) to match them too.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Paragraph border: top/left/bottom/right, each just reserving 5pt of
# space (no visible rule) - matches the pBdr already on paragraphs 3/4.
$b = $p1.Range.ParagraphFormat.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Left indent 120 twips (6pt) -> 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Retarget the ID marker text and drop the now-redundant trailing
# space run that used to trail it.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5316_topic_8__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5316_506__ID**", 2)
